$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 532.1786
$ws.Range("I53").Value = 216
$ws.Range("J53").Value = 1101.3
$ws.Range("K53").Value = 216
$ws.Range("L53").Value = 1101.3
$ws.Range("M53").Value = 421
$ws.Range("N53").Value = -2375.3

$ws.Range("H64").Value = 4288
$ws.Range("I64").Value = 4454
$ws.Range("J64").Value = 3956
$ws.Range("K64").Value = 4454
$ws.Range("L64").Value = 3956
$ws.Range("M64").Value = -4206
$ws.Range("N64").Value = -4452

$ws.Range("H67").Value = 4288
$ws.Range("I67").Value = 4454
$ws.Range("J67").Value = 3956
$ws.Range("K67").Value = 4454
$ws.Range("L67").Value = 3956
$ws.Range("M67").Value = -3596
$ws.Range("N67").Value = -5672

$ws.Range("H76").Value = 3106.976
$ws.Range("I76").Value = 3064.5
$ws.Range("K76").Value = 3064.5
$ws.Range("M76").Value = -2749.5

$ws.Range("H79").Value = 3106.976
$ws.Range("I79").Value = 3064.5
$ws.Range("K79").Value = 3064.5
$ws.Range("M79").Value = -1972.5

$ws.Range("H98").Value = 2673.7354
$ws.Range("I98").Value = 2831.1936
$ws.Range("K98").Value = 2831.1936
$ws.Range("M98").Value = -1333.1936

$ws.Range("H122").Value = 2673.7354
$ws.Range("I122").Value = 2831.1936
$ws.Range("K122").Value = 8493.5808
$ws.Range("M122").Value = -6043.5808

$ws.Range("H129").Value = 450.1
$ws.Range("I129").Value = 298.875
$ws.Range("J129").Value = 1055
$ws.Range("K129").Value = 896.625
$ws.Range("L129").Value = 3165
$ws.Range("M129").Value = 4103.375
$ws.Range("N129").Value = -13165

$ws.Range("H135").Value = 11044503
$ws.Range("I135").Value = 2421.225
$ws.Range("J135").Value = 66254910
$ws.Range("K135").Value = 21791.025
$ws.Range("L135").Value = 596294190
$ws.Range("M135").Value = -19256.025
$ws.Range("N135").Value = -596299260

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7769.87
$ws.Range("I32").Value = 5727.9663
$ws.Range("J32").Value = 24290.727
$ws.Range("K32").Value = 5727.9663
$ws.Range("L32").Value = 24290.727
$ws.Range("M32").Value = -5440.9663
$ws.Range("N32").Value = -24864.727

$ws.Range("H61").Value = 4834991.5
$ws.Range("I61").Value = 5852492.5
$ws.Range("J61").Value = 1860
$ws.Range("K61").Value = 5852492.5
$ws.Range("L61").Value = 1860
$ws.Range("M61").Value = -5852280.5
$ws.Range("N61").Value = -2284

$ws.Range("H63").Value = 3373.4
$ws.Range("I63").Value = 2841.75
$ws.Range("K63").Value = 2841.75
$ws.Range("M63").Value = -2155.75

$ws.Range("H66").Value = 3373.4
$ws.Range("I66").Value = 2841.75
$ws.Range("K66").Value = 14208.75
$ws.Range("M66").Value = -10776.75

$ws.Range("H88").Value = 2491.4285
$ws.Range("I88").Value = 2725
$ws.Range("J88").Value = 2180
$ws.Range("K88").Value = 2725
$ws.Range("L88").Value = 2180
$ws.Range("M88").Value = -2319
$ws.Range("N88").Value = -2992

$ws.Range("H91").Value = 2491.4285
$ws.Range("I91").Value = 2725
$ws.Range("J91").Value = 2180
$ws.Range("K91").Value = 2725
$ws.Range("L91").Value = 2180
$ws.Range("M91").Value = -1321
$ws.Range("N91").Value = -4988

$ws.Range("H110").Value = 1495.8125
$ws.Range("I110").Value = 1334.7778
$ws.Range("J110").Value = 1702.8572
$ws.Range("K110").Value = 1334.7778
$ws.Range("L110").Value = 1702.8572
$ws.Range("M110").Value = 710.2221999999999
$ws.Range("N110").Value = -5792.8572

$ws.Range("H136").Value = 4834991.5
$ws.Range("I136").Value = 5852492.5
$ws.Range("J136").Value = 1860
$ws.Range("K136").Value = 17557477.5
$ws.Range("L136").Value = 5580
$ws.Range("M136").Value = -17554927.5
$ws.Range("N136").Value = -10680

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1556.2222
$ws.Range("I86").Value = 1480.762
$ws.Range("K86").Value = 1480.762
$ws.Range("M86").Value = -357.7619999999999

$ws.Range("H89").Value = 1556.2222
$ws.Range("I89").Value = 1480.762
$ws.Range("K89").Value = 7403.809999999999
$ws.Range("M89").Value = -1787.809999999999

$ws.Range("H105").Value = 1724.9286
$ws.Range("I105").Value = 1695.3636
$ws.Range("K105").Value = 1695.3636
$ws.Range("M105").Value = 51.63640000000009

$ws.Range("H134").Value = 8745207
$ws.Range("I134").Value = 13408448
$ws.Range("J134").Value = 1630.5
$ws.Range("K134").Value = 40225344
$ws.Range("L134").Value = 4891.5
$ws.Range("M134").Value = -40222809
$ws.Range("N134").Value = -9961.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1417
$ws.Range("I31").Value = 999.3946999999999
$ws.Range("J31").Value = 2637.6924
$ws.Range("K31").Value = 999.3946999999999
$ws.Range("L31").Value = 2637.6924
$ws.Range("M31").Value = -704.3946999999999
$ws.Range("N31").Value = -3227.6924

$ws.Range("H34").Value = 1417
$ws.Range("I34").Value = 999.3946999999999
$ws.Range("J34").Value = 2637.6924
$ws.Range("K34").Value = 999.3946999999999
$ws.Range("L34").Value = 2637.6924
$ws.Range("M34").Value = -797.3946999999999
$ws.Range("N34").Value = -3041.6924

$ws.Range("H52").Value = 47086.668
$ws.Range("J52").Value = 47086.668
$ws.Range("L52").Value = 47086.668
$ws.Range("N52").Value = -47674.668

$ws.Range("H62").Value = 2742.7222
$ws.Range("I62").Value = 2407.182
$ws.Range("J62").Value = 3270
$ws.Range("K62").Value = 2407.182
$ws.Range("L62").Value = 3270
$ws.Range("M62").Value = -1783.182
$ws.Range("N62").Value = -4518

$ws.Range("H65").Value = 2742.7222
$ws.Range("I65").Value = 2407.182
$ws.Range("J65").Value = 3270
$ws.Range("K65").Value = 12035.91
$ws.Range("L65").Value = 16350
$ws.Range("M65").Value = -8915.91
$ws.Range("N65").Value = -22590

$ws.Range("H141").Value = 31792
$ws.Range("J141").Value = 31792
$ws.Range("L141").Value = 31792
$ws.Range("N141").Value = -42152

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 1300
$ws.Range("I81").Value = 500
$ws.Range("K81").Value = 1500
$ws.Range("M81").Value = -377

$ws.Range("H84").Value = 1300
$ws.Range("I84").Value = 500
$ws.Range("K84").Value = 4500
$ws.Range("M84").Value = 1116

$ws.Range("H131").Value = 2104.8408
$ws.Range("I131").Value = 6456
$ws.Range("J131").Value = 1547
$ws.Range("K131").Value = 19368
$ws.Range("L131").Value = 4641
$ws.Range("M131").Value = -14328
$ws.Range("N131").Value = -14721

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8140.963
$ws.Range("I70").Value = 10819.8
$ws.Range("J70").Value = 4792.4165
$ws.Range("K70").Value = 10819.8
$ws.Range("L70").Value = 4792.4165
$ws.Range("M70").Value = -10549.8
$ws.Range("N70").Value = -5332.4165

$ws.Range("H73").Value = 8140.963
$ws.Range("I73").Value = 10819.8
$ws.Range("J73").Value = 4792.4165
$ws.Range("K73").Value = 10819.8
$ws.Range("L73").Value = 4792.4165
$ws.Range("M73").Value = -9883.799999999999
$ws.Range("N73").Value = -6664.4165

$ws.Range("H80").Value = 12535.5
$ws.Range("I80").Value = 2888.125
$ws.Range("J80").Value = 51125
$ws.Range("K80").Value = 2888.125
$ws.Range("L80").Value = 51125
$ws.Range("M80").Value = -1890.125
$ws.Range("N80").Value = -53121

$ws.Range("H83").Value = 12535.5
$ws.Range("I83").Value = 2888.125
$ws.Range("J83").Value = 51125
$ws.Range("K83").Value = 14440.625
$ws.Range("L83").Value = 255625
$ws.Range("M83").Value = -9448.625
$ws.Range("N83").Value = -265609

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3862.074
$ws.Range("I132").Value = 4600.85
$ws.Range("J132").Value = 1751.2858
$ws.Range("K132").Value = 13802.55
$ws.Range("L132").Value = 5253.857400000001
$ws.Range("M132").Value = -11272.55
$ws.Range("N132").Value = -10313.8574

$ws.Range("H136").Value = 11144.833
$ws.Range("I136").Value = 13361.173
$ws.Range("J136").Value = 1962.8572
$ws.Range("K136").Value = 40083.519
$ws.Range("L136").Value = 5888.571599999999
$ws.Range("M136").Value = -37533.519
$ws.Range("N136").Value = -10988.5716
